$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "B" (Jun_13) column, pushing
# the old B -> D and old C -> E. This mirrors selecting columns B:C and
# choosing Insert in the Excel UI.
$ws.Columns("B:C").Insert()

# New header row values for the freshly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# The inserted columns are populated the same way column B originally was:
# every data row gets the "UN" placeholder.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Match the width already used by the (now shifted) data column so the two
# new columns line up visually with the rest of the sheet.
$ws.Columns("C").ColumnWidth = 7.1667
$ws.Columns("D").ColumnWidth = 7.1667
$ws.Columns("E").ColumnWidth = 7.1667
